# Normalize casing of categorical values in the shark-attack style dataset.
# Columns E (type), G (sex), and J (fatal) sometimes hold the literal values
# Unknown / Unprovoked / Provoked / Watercraft which should be lower-cased.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 1) { $lastRow = 1 }

$targets = @("Unknown", "Unprovoked", "Provoked", "Watercraft")
$cols = @(5, 7, 10)  # E, G, J

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $targets -contains $val) {
            $cell.Value2 = $val.ToString().ToLower()
        }
    }
}
